# Atualiza testes da planilha
# Target sheet: "Produto" (second worksheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Update ordem do teste for rows 11 and 12 (TC009 / TC010)
$ws.Range("B11").Value = 9
$ws.Range("B12").Value = 10

# Insert 6 new rows (13-18) for the new test cases TC011-TC016, pushing the
# blank separator row and the "Observações adicionais" block further down.
$ws.Rows("13:18").Insert()

# Copy the formatting (borders/fonts/fill) of an existing fully-bordered data
# row (row 11, all columns styled) onto the newly inserted rows so every cell
# A:I gets the same "data row" style used elsewhere in the table.
$ws.Range("A11:I11").Copy()
$ws.Range("A13:I18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 13 - TC011
$ws.Range("A13").Value = "TC011"
$ws.Range("B13").Value = 11
$ws.Range("C13").Value = "N/A"
$ws.Range("D13").Value = "N/A"
$ws.Range("E13").Value = "N/A"
$ws.Range("F13").Value = "N/A"
$ws.Range("G13").Value = "N/A"
$ws.Range("H13").Value = "Deve abrir o modal de cadastro de produto ao clicar no botão de criar "
$ws.Range("I13").Value = "Falha"

# Row 14 - TC012
$ws.Range("A14").Value = "TC012"
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "N/A"
$ws.Range("G14").Value = "N/A"
$ws.Range("H14").Value = "Deve abrir o modal de cadastro de produto ao clicar no botão de editar "
$ws.Range("I14").Value = "Falha"

# Row 15 - TC013
$ws.Range("A15").Value = "TC013"
$ws.Range("B15").Value = 13
$ws.Range("C15").Value = "N/A"
$ws.Range("D15").Value = "N/A"
$ws.Range("E15").Value = "N/A"
$ws.Range("F15").Value = "N/A"
$ws.Range("G15").Value = "N/A"
$ws.Range("H15").Value = "Deve excluir o produto ao clicar no botão de excluir"
$ws.Range("I15").Value = "Falha"

# Row 16 - TC014
$ws.Range("A16").Value = "TC014"
$ws.Range("B16").Value = 14
$ws.Range("C16").Value = "N/A"
$ws.Range("D16").Value = "N/A"
$ws.Range("E16").Value = "N/A"
$ws.Range("F16").Value = "N/A"
$ws.Range("G16").Value = "N/A"
$ws.Range("H16").Value = "Deve voltar para a página de login ao clicar no botão de voltar "
$ws.Range("I16").Value = "Falha"

# Row 17 - TC015
$ws.Range("A17").Value = "TC015"
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = "N/A"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = "N/A"
$ws.Range("F17").Value = "N/A"
$ws.Range("G17").Value = "N/A"
$ws.Range("H17").Value = "Deve emitir uma mensagem de erro ao tentar cadastrar um produto com mesmo código de um já existente"
$ws.Range("I17").Value = "Falha"

# Row 18 - TC016
$ws.Range("A18").Value = "TC016"
$ws.Range("B18").Value = 16
$ws.Range("C18").Value = "N/A"
$ws.Range("D18").Value = "N/A"
$ws.Range("E18").Value = "N/A"
$ws.Range("F18").Value = "N/A"
$ws.Range("G18").Value = "N/A"
$ws.Range("H18").Value = "Deve fechar o modal ao clicar no botão de sair "
$ws.Range("I18").Value = "Sucesso"

# Update the "Observações adicionais" note that used to live in row 16 (now
# shifted to row 22) with the new wording.
$ws.Range("A22").Value = "1. Equipe de desenvolvimento deve padronizar a lingua da página;"
